# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.88
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AG3").Value = 8.5
$ws.Range("AH3").Value = 19
$ws.Range("AN3").Value = 3.75
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 101

# Row 8 changes
$ws.Range("G8").Value = 1.29
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 8.5
$ws.Range("L8").Value = 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("S8").Value = 1.29
$ws.Range("T8").Value = 3.5
$ws.Range("W8").Value = 8
$ws.Range("AD8").Value = 9.5
$ws.Range("AE8").Value = 21
$ws.Range("AT8").Value = 3.5
$ws.Range("AW8").Value = 9.5
